$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the Price/Volume columns so numeric-looking
# strings (e.g. "0.9972", "241.17") are not auto-coerced to numbers,
# matching the original inline-string cell type.
$range = $ws.Range("D2:E51")
$range.NumberFormat = "@"

$ws.Range("D2").Value = '29.097.57'
$ws.Range("E2").Value = '  -0.07%  '
$ws.Range("D3").Value = '1.821.69'
$ws.Range("E3").Value = '  -0.62%  '
$ws.Range("D4").Value = '0.9972'
$ws.Range("E4").Value = '  -0.32%  '
$ws.Range("D5").Value = '241.17'
$ws.Range("E5").Value = '  -0.93%  '
$ws.Range("D6").Value = '0.6148'
$ws.Range("E6").Value = '  -2.13%  '
$ws.Range("D7").Value = '0.9989'
$ws.Range("E7").Value = '  -0.37%  '
$ws.Range("D8").Value = '0.07324'
$ws.Range("E8").Value = '  -2.20%  '
$ws.Range("D9").Value = '0.2887'
$ws.Range("E9").Value = '  -1.28%  '
$ws.Range("D10").Value = '22.93'
$ws.Range("E10").Value = '  -1.14%  '
$ws.Range("D11").Value = '0.07647'
$ws.Range("E11").Value = '  -0.42%  '
$ws.Range("D12").Value = '1.823.85'
$ws.Range("E12").Value = '  -0.21%  '
$ws.Range("D13").Value = '4.944'
$ws.Range("E13").Value = '  -1.30%  '
$ws.Range("D14").Value = '0.6587'
$ws.Range("E14").Value = '  -1.39%  '
$ws.Range("D15").Value = '81.73'
$ws.Range("E15").Value = '  -1.26%  '
$ws.Range("D16").Value = '0.000008916'
$ws.Range("E16").Value = '  -5.27%  '
$ws.Range("D17").Value = '5.822'
$ws.Range("E17").Value = '  -2.67%  '
$ws.Range("D18").Value = '29.063.76'
$ws.Range("E18").Value = '  -0.20%  '
$ws.Range("D19").Value = '2.068.84'
$ws.Range("E19").Value = '  -0.54%  '
$ws.Range("D20").Value = '236.75'
$ws.Range("E20").Value = '  +6.11%  '
$ws.Range("D21").Value = '12.42'
$ws.Range("E21").Value = '  -1.31%  '
$ws.Range("D22").Value = '0.9985'
$ws.Range("E22").Value = '  -0.52%  '
$ws.Range("D23").Value = '7.113'
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("D24").Value = '0.9980'
$ws.Range("E24").Value = '  -0.37%  '
$ws.Range("D25").Value = '157.33'
$ws.Range("E25").Value = '  -1.66%  '
$ws.Range("D26").Value = '0.1411'
$ws.Range("E26").Value = '  +1.42%  '
$ws.Range("D27").Value = '8.413'
$ws.Range("E27").Value = '  -0.92%  '
$ws.Range("D28").Value = '17.59'
$ws.Range("E28").Value = '  -1.55%  '
$ws.Range("D29").Value = '1.483'
$ws.Range("E29").Value = '  -0.94%  '
$ws.Range("D30").Value = '0.05545'
$ws.Range("E30").Value = '  -3.22%  '
$ws.Range("D31").Value = '4.083'
$ws.Range("E31").Value = '  -0.03%  '
$ws.Range("D32").Value = '4.084'
$ws.Range("E32").Value = '  -1.58%  '
$ws.Range("D33").Value = '1.205'
$ws.Range("E33").Value = '  -0.29%  '
$ws.Range("D34").Value = '1.819'
$ws.Range("E34").Value = '  -0.61%  '
$ws.Range("D35").Value = '0.7327'
$ws.Range("E35").Value = '  -0.98%  '
$ws.Range("D36").Value = '1.130'
$ws.Range("E36").Value = '  -0.81%  '
$ws.Range("D37").Value = '2.604'
$ws.Range("E37").Value = '  -2.60%  '
$ws.Range("E38").Value = '  +1.94%  '
$ws.Range("D39").Value = '1.204.99'
$ws.Range("E39").Value = '  -0.80%  '
$ws.Range("D40").Value = '0.01753'
$ws.Range("E40").Value = '  -1.38%  '
$ws.Range("D41").Value = '6.330'
$ws.Range("E41").Value = '  -2.73%  '
$ws.Range("D42").Value = '0.8976'
$ws.Range("E42").Value = '  +1.01%  '
$ws.Range("D43").Value = '0.9988'
$ws.Range("E43").Value = '  -0.41%  '
$ws.Range("D44").Value = '101.09'
$ws.Range("E44").Value = '  -0.99%  '
$ws.Range("D45").Value = '1.972.22'
$ws.Range("E45").Value = '  -0.38%  '
$ws.Range("D46").Value = '64.46'
$ws.Range("E46").Value = '  -1.61%  '
$ws.Range("D47").Value = '0.5077'
$ws.Range("E47").Value = '  -0.27%  '
$ws.Range("D48").Value = '0.00000000119'
$ws.Range("E48").Value = '  -3.88%  '
$ws.Range("D49").Value = '0.3993'
$ws.Range("E49").Value = '  -1.69%  '
$ws.Range("D50").Value = '8.982'
$ws.Range("E50").Value = '  +0.41%  '
$ws.Range("D51").Value = '0.05747'
$ws.Range("E51").Value = '  -1.29%  '

# Restore the default (unstyled) cell style so no stray number-format
# style survives on these cells, matching the target which has no
# style attribute on D/E data cells.
$range.Style = "Normal"

